$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "52.110.35"
Set-TextValue "E2" "  +0.78%  "
Set-TextValue "D3" "2.897.51"
Set-TextValue "E3" "  +3.57%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "351.39"
Set-TextValue "E5" "  -0.51%  "
Set-TextValue "D6" "112.57"
Set-TextValue "E6" "  +0.98%  "
Set-TextValue "D7" "0.556"
Set-TextValue "E7" "  -0.15%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.620"
Set-TextValue "E9" "  -0.46%  "
Set-TextValue "D10" "39.72"
Set-TextValue "E10" "  -1.30%  "
Set-TextValue "E11" "  +0.69%  "
Set-TextValue "D12" "0.0859"
Set-TextValue "E12" "  +2.78%  "
Set-TextValue "D13" "19.72"
Set-TextValue "E13" "  -0.83%  "
Set-TextValue "E14" "  -0.94%  "
Set-TextValue "D15" "3.354.73"
Set-TextValue "E15" "  +3.64%  "
Set-TextValue "D16" "2.904.10"
Set-TextValue "E16" "  +3.68%  "
Set-TextValue "D17" "0.982"
Set-TextValue "E17" "  +4.13%  "
Set-TextValue "D18" "52.185.66"
Set-TextValue "E18" "  +0.96%  "
Set-TextValue "D20" "7.59"
Set-TextValue "E20" "  -0.18%  "
Set-TextValue "D21" "13.89"
Set-TextValue "E21" "  +2.37%  "
Set-TextValue "D22" "0.0₃0973"
Set-TextValue "E22" "  +0.30%  "
Set-TextValue "D23" "70.79"
Set-TextValue "E23" "  +0.75%  "
Set-TextValue "D24" "268.41"
Set-TextValue "E24" "  +0.48%  "
Set-TextValue "E25" "  +1.05%  "
Set-TextValue "D26" "0.179"
Set-TextValue "E26" "  +12.41%  "
Set-TextValue "E27" "  +2.15%  "
Set-TextValue "E28" "  -0.05%  "
Set-TextValue "D29" "10.61"
Set-TextValue "E29" "  +2.48%  "
Set-TextValue "D30" "0.103"
Set-TextValue "E30" "  +15.51%  "
Set-TextValue "D31" "6.60"
Set-TextValue "E31" "  +7.75%  "
Set-TextValue "D32" "37.32"
Set-TextValue "E32" "  -4.08%  "
Set-TextValue "E33" "  -0.70%  "
Set-TextValue "D34" "6.18"
Set-TextValue "E34" "  +11.62%  "
Set-TextValue "D35" "52.89"
Set-TextValue "E35" "  +0.82%  "
Set-TextValue "D36" "0.0450"
Set-TextValue "E36" "  -0.53%  "
Set-TextValue "D37" "0.999"
Set-TextValue "E37" "  -0.10%  "
Set-TextValue "D38" "3.30"
Set-TextValue "E38" "  +4.41%  "
Set-TextValue "D39" "18.81"
Set-TextValue "E39" "  -0.02%  "
Set-TextValue "E40" "  +1.26%  "
Set-TextValue "E41" "  +8.32%  "
Set-TextValue "D42" "0.116"
Set-TextValue "E42" "  +1.19%  "
Set-TextValue "D43" "22.98"
Set-TextValue "E43" "  +4.82%  "
Set-TextValue "B44" "Monero"
Set-TextValue "C44" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D44" "119.54"
Set-TextValue "E44" "  -0.59%  "
Set-TextValue "B45" "ApeXProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D45" "2.59"
Set-TextValue "E45" "  +5.06%  "
Set-TextValue "E46" "  -1.88%  "
Set-TextValue "B47" "NEARProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D47" "3.49"
Set-TextValue "E47" "  +1.53%  "
Set-TextValue "B48" "Maker"
Set-TextValue "C48" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D48" "2.167.36"
Set-TextValue "E48" "  +2.97%  "
Set-TextValue "D49" "0.263"
Set-TextValue "E49" "  +19.96%  "
Set-TextValue "D50" "0.0342"
Set-TextValue "E50" "  +10.21%  "
Set-TextValue "D51" "0.948"
Set-TextValue "E51" "  -0.65%  "
